# Update first page of excel sheet
#
# The workbook has two sheets:
#   1) "Critical Path Analysis"  (was not the active tab)
#   2) "Risky Path Analysis"     (was the active tab, scrolled to column B,
#                                 zoomed to 106%)
#
# This script:
#   - Fills in the previously-empty ESD/LED (columns D/E) values on the
#     "Critical Path Analysis" sheet.
#   - Makes "Critical Path Analysis" the active/selected sheet, scrolled
#     back to A1, zoomed to 150%, with E3 selected.
#   - Leaves "Risky Path Analysis" selection as-is (N37) but it is no
#     longer the active tab, and its horizontal scroll resets off of B1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Critical Path Analysis")

# --- Fill in the Earliest Start Date (D) / Latest End Date (E) columns ---
$ws1.Range("E2").Value  = 1

$ws1.Range("D3").Value  = 1
$ws1.Range("E3").Value  = 3

$ws1.Range("D4").Value  = 1
$ws1.Range("E4").Value  = 3

$ws1.Range("D5").Value  = 3
$ws1.Range("E5").Value  = 7

$ws1.Range("D6").Value  = 3
$ws1.Range("E6").Value  = 6

$ws1.Range("D7").Value  = 7
$ws1.Range("E7").Value  = 13

$ws1.Range("D8").Value  = 6
$ws1.Range("E8").Value  = 10

$ws1.Range("D9").Value  = 6
$ws1.Range("E9").Value  = 10

$ws1.Range("D10").Value = 10
$ws1.Range("E10").Value = 13

$ws1.Range("D11").Value = 10
$ws1.Range("E11").Value = 13

$ws1.Range("D12").Value = 13
$ws1.Range("E12").Value = 16

# --- Switch the active/selected sheet to "Critical Path Analysis" ---
$ws1.Activate() | Out-Null
$ws1.Range("E3").Select() | Out-Null
$excel.ActiveWindow.Zoom = 150
